$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New WRESBAL data (date serial, value) replacing the old observation window.
$data = @(
    @(44475, 4130.05),
    @(44482, 4194.503),
    @(44489, 4134.714),
    @(44496, 4077.17),
    @(44503, 4086.336),
    @(44510, 4169.041),
    @(44517, 4209.802),
    @(44524, 4162.81),
    @(44531, 4224.865),
    @(44538, 4275.806),
    @(44545, 4251.431),
    @(44552, 4185.009),
    @(44559, 4116.196),
    @(44566, 3859.731),
    @(44573, 3973.554),
    @(44580, 3932.521),
    @(44587, 3796.174),
    @(44594, 3789.714),
    @(44601, 3771.911),
    @(44608, 3799.973),
    @(44615, 3769.869),
    @(44622, 3844.716),
    @(44629, 3959.659),
    @(44636, 3961.721),
    @(44643, 3801.095),
    @(44650, 3829.68),
    @(44657, 3834.183),
    @(44664, 3823.644),
    @(44671, 3652.101),
    @(44678, 3331.982),
    @(44685, 3296.085),
    @(44692, 3288.552),
    @(44699, 3356.377),
    @(44706, 3306.18),
    @(44713, 3322.03),
    @(44720, 3339.405),
    @(44727, 3283.954),
    @(44734, 3147.556),
    @(44741, 3130.087),
    @(44748, 3190.627),
    @(44755, 3273.436),
    @(44762, 3284.324),
    @(44769, 3256.875),
    @(44776, 3259.011),
    @(44783, 3348.147),
    @(44790, 3346.75),
    @(44797, 3307.365),
    @(44804, 3248.654),
    @(44811, 3285.498),
    @(44818, 3234.928),
    @(44825, 3104.032),
    @(44832, 2960.598),
    @(44839, 2972.874),
    @(44846, 3087.912),
    @(44853, 3092.044),
    @(44860, 3042.838),
    @(44867, 3062.727),
    @(44874, 3067.888),
    @(44881, 3128.277),
    @(44888, 3160.73),
    @(44895, 3168.995),
    @(44902, 3199.048),
    @(44909, 3188.634),
    @(44916, 3088.621),
    @(44923, 3017.889),
    @(44930, 2830.144),
    @(44937, 3081.237),
    @(44944, 3147.332),
    @(44951, 3031.889),
    @(44958, 3015.506),
    @(44965, 3037.584),
    @(44972, 3041.038),
    @(44979, 3016.577),
    @(44986, 2998.303),
    @(44993, 2999.72),
    @(45000, 3251.482),
    @(45007, 3425.075),
    @(45014, 3437.592),
    @(45021, 3292.058),
    @(45028, 3413.989),
    @(45035, 3303.152),
    @(45042, 3135.655),
    @(45049, 3089.706),
    @(45056, 3196.429),
    @(45063, 3248.878),
    @(45070, 3240.172),
    @(45077, 3309.447),
    @(45084, 3349.751),
    @(45091, 3325.261),
    @(45098, 3251.558),
    @(45105, 3171.567),
    @(45112, 3114.989)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Range("A" + $row).Value = $data[$i][0]
    $ws.Range("B" + $row).Value = $data[$i][1]
}

$newLastRow = $data.Count + 1
$oldUsedRows = $ws.UsedRange.Rows.Count
if ($oldUsedRows -gt $newLastRow) {
    $deleteRange = ($newLastRow + 1).ToString() + ":" + $oldUsedRows.ToString()
    $ws.Rows($deleteRange).Delete()
}

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B3").Value = "2023-07-09"
$ws2.Range("B4").Value = "2023-07-09"
$ws2.Range("B7").Value = "2023-07-05"
$ws2.Range("B14").Value = "2023-07-06 15:35:25-05"
$ws2.Range("B15").Value = 78
